$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "zika",
  "flu",
  "influenza",
  "hiv",
  "vis",
  "malaria",
  "hpv",
  "ebola",
  "diabetes",
  "hepatitis b",
  "tuberculosis",
  "obesity",
  "vaccines",
  "chlamydia",
  "immunization schedule 2017",
  "zika virus",
  "influenza 2017",
  "measles",
  "mumps",
  "rabies"
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

Write-Host "done"
